$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 17.83579802019262
$ws.Range("F2").Value = 0.5702543166250014
$ws.Range("K2").Value = 0.5795058394867494
$ws.Range("E3").Value = 6.949933224568027
$ws.Range("F3").Value = 0.3567451747526734
$ws.Range("K3").Value = 0.541994255313475
$ws.Range("E4").Value = 20.0908989192975
$ws.Range("F4").Value = 0.7604596895199529
$ws.Range("K4").Value = 0.2384700568738055
$ws.Range("E5").Value = 22.36722787308386
$ws.Range("F5").Value = 0.8466209118920863
$ws.Range("K5").Value = 0.369987492810365
$ws.Range("E6").Value = 20.0908989192975
$ws.Range("F6").Value = 0.7604596895199529
$ws.Range("K6").Value = 0.2384700568738055
$ws.Range("E7").Value = 22.36722787308386
$ws.Range("F7").Value = 0.8466209118920863
$ws.Range("K7").Value = 0.369987492810365
$ws.Range("E8").Value = 9.067515113986417
$ws.Range("F8").Value = 0.5224032507173562
$ws.Range("K8").Value = 0.3531987790792718
$ws.Range("E9").Value = 2.848304889328497
$ws.Range("F9").Value = 0.2804547627048343
$ws.Range("K9").Value = 0.4455508080542069
$ws.Range("E10").Value = 12.93506843424415
$ws.Range("F10").Value = 0.48982208828641
$ws.Range("K10").Value = 0.994277316800055
$ws.Range("E11").Value = 11.21102148593617
$ws.Range("F11").Value = 0.3491636681036706
$ws.Range("K11").Value = 0.7838758923213196
$ws.Range("E12").Value = 12.50281531047952
$ws.Range("F12").Value = 0.7675496598637382
$ws.Range("K12").Value = 0.398698891381401
$ws.Range("E13").Value = 7.104774546832009
$ws.Range("F13").Value = 0.4870637756620864
$ws.Range("K13").Value = 0.434694428313129
$ws.Range("E14").Value = 2.765203229139217
$ws.Range("F14").Value = 0.2979157940111148
$ws.Range("K14").Value = 0.2591662631864758
$ws.Range("E15").Value = 2.680320957381447
$ws.Range("F15").Value = 0.2730483328392531
$ws.Range("K15").Value = 0.2718126216402545
$ws.Range("E16").Value = 12.54392383586099
$ws.Range("F16").Value = 0.3906764842469256
$ws.Range("K16").Value = 0.9856959691366205
$ws.Range("E17").Value = 24.27650286081777
$ws.Range("F17").Value = 0.9085859753463218
$ws.Range("K17").Value = 0.8715344932517999
$ws.Range("E18").Value = 10.43527287151567
$ws.Range("F18").Value = 0.6142701292199257
$ws.Range("K18").Value = 0.4607920932405188
$ws.Range("E19").Value = 3.974395028628913
$ws.Range("F19").Value = 0.4546606527283644
$ws.Range("K19").Value = 0.2710238717441542
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("K20").Value = "inf"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = -0.2904283587044502
$ws.Range("K21").Value = "inf"
$ws.Range("E22").Value = 12.50281576818371
$ws.Range("F22").Value = 0.541645262347088
$ws.Range("K22").Value = 0.3815100007340863
$ws.Range("E23").Value = 2.516212267075455
$ws.Range("F23").Value = 0.3482726110306784
$ws.Range("K23").Value = 0.1955980418471334
$ws.Range("E24").Value = 3.974395028628931
$ws.Range("F24").Value = 0.3811634430807067
$ws.Range("K24").Value = 0.2710276536371574
$ws.Range("E25").Value = 2.961217231488562
$ws.Range("F25").Value = 0.3064665867744434
$ws.Range("K25").Value = 0.4471919301970728
$ws.Range("E26").Value = 9.737650172584237
$ws.Range("F26").Value = 0.7125676451573347
$ws.Range("K26").Value = 0.4687496917930938
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("K27").Value = "inf"
$ws.Range("E28").Value = 6.774435839264163
$ws.Range("F28").Value = 0.2964700727530642
$ws.Range("K28").Value = 0.441454738466561
$ws.Range("E29").Value = 0.1333366471518599
$ws.Range("F29").Value = 0.2771660039182953
$ws.Range("K29").Value = 0.4431278979244214
$ws.Range("N29").Value = "c_sources:EX_glc__D_e;BIOMASS_AERO_SC_RBA:(0.009615384615,0.14);EX_glc__D_e:(-2.20113372715757,1000);EX_nh4_e:(-0.181824845165991,1000);EX_o2_e:(-26.25403530509562,1000)"
$ws.Range("E30").Value = -0.0000008946889123367239
$ws.Range("F30").Value = 0.00000845935646720553
$ws.Range("N30").Value = "c_sources:EX_xyl__D_e,EX_glc__D_e;EX_xyl__D_e:(-0.0266525306904169,-0.0266525306904169);EX_glc__D_e:(-0.0518243652313661,-0.0518243652313661);BIOMASS_AERO_SC_RBA:(0.0138888888888888,0.0138888888888889);ATPM_c:(0,1000);EX_4abz_e:(-1000,1000);EX_fe3_e:(-1000,1000);EX_fol_e:(-1000,1000);EX_inost_e:(-1000,1000);EX_nac_e:(-1000,1000);EX_pnto__R_e:(-1000,1000);EX_pydxn_e:(-1000,1000);EX_ribflv_e:(-1000,1000);EX_thm_e:(-1000,1000);EX_o2_e:(-26.25403530509562,1000)"
